$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 768 ("凍結路面を運転中、彼女は嫌な場面に居合わせた" post) was removed from the sheet.
# Deleting the entire row shifts all subsequent rows up by one.
$ws.Rows(768).Delete()
